$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 58, shifting existing rows 58-64 down to 59-65
$ws.Rows.Item(58).Insert()

# Populate the new row 58 with the new weekly data record
$ws.Range("A58").Value = 10
$ws.Range("B58").Value = "Vega Modelo de Temuco"
$ws.Range("C58").Value = "La Araucanía"
$ws.Range("D58").Value = 44301
$ws.Range("E58").Value = 9
$ws.Range("F58").Value = "Fruta"
$ws.Range("G58").Value = 100108
$ws.Range("H58").Value = "Tropicales y subtropicales"
$ws.Range("I58").Value = 100108004
$ws.Range("J58").Value = "Papaya"
$ws.Range("K58").Value = "Cultivar IV Región"
$ws.Range("L58").Value = "Primera"
$ws.Range("M58").Value = 55
$ws.Range("N58").Value = 21000
$ws.Range("O58").Value = 21000
$ws.Range("P58").Value = 21000
$ws.Range("Q58").Value = "`$/caja 15 kilos granel"
$ws.Range("R58").Value = "Provincia del Elquí"
$ws.Range("S58").Value = 1400
$ws.Range("T58").Value = 15

# Copy the date cell style (D column) from the row above to keep formatting consistent
$ws.Range("D57").Copy()
$ws.Range("D58").PasteSpecial(-4122)
